# Update ticker data: append new rows 47-58 to Sheet1, column A
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTickers = @(
    "DOT-USD",
    "AVAX-USD",
    "SHIB-USD",
    "BUSD-USD",
    "UNI-USD",
    "ICP-USD",
    "LDO-USD",
    "GRT-USD",
    "AAVE-USD",
    "FRAX-USD",
    "RETH-USD",
    "EGLD-USD"
)

$startRow = 47
for ($i = 0; $i -lt $newTickers.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newTickers[$i]
}
